$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.581.72"
$ws.Range("E2").Value = "  -1.13%  "

$ws.Range("D3").Value = "1.847.50"
$ws.Range("E3").Value = "  -2.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -1.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.94"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4602"
$ws.Range("E7").Value = "  -1.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3884"
$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.99"
$ws.Range("E9").Value = "  -2.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07937"
$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.56"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").Value = "1.859.35"
$ws.Range("E13").Value = "  -1.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.975"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.165"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.43"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06675"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001037"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.20"
$ws.Range("E20").Value = "  +1.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("D22").Value = "27.620.63"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.407"
$ws.Range("E23").Value = "  -1.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").Value = "2.079.70"
$ws.Range("E26").Value = "  -1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.19"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("E28").Value = "  -2.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("E29").Value = "  +3.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.445"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.75"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9767"
$ws.Range("E32").Value = "  +1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09404"
$ws.Range("E33").Value = "  -0.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.303"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.339"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02232"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06018"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.302"
$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.185"
$ws.Range("E40").Value = "  -2.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5924"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1867"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.35"
$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5592"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.16"
$ws.Range("E46").Value = "  +0.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.912"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06703"
$ws.Range("E48").Value = "  -2.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.68"
$ws.Range("E49").Value = "  -2.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.051"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("E51").Value = "  -1.24%  "
